$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be created in this order to match target string table:
# cadmiumoleate, topselenide, Cadmium Oleate, Triphenylphosphine selenide
$ws.Range("A7").Value = "cadmiumoleate"
$ws.Range("A8").Value = "topselenide"
$ws.Range("B7").Value = "Cadmium Oleate"
$ws.Range("B8").Value = "Triphenylphosphine selenide"

$ws.Range("C7").Value = 675.3
$ws.Range("D7").Value = "Negligible"

$ws.Range("C8").Value = 341.3
$ws.Range("D8").Value = "Negligible"

# Selection moved to H28 per diff
$ws.Range("H28").Select()
